$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds plain-looking decimal numbers stored as
# TEXT in the source data (e.g. "2.60", "33.03") so that formatting such
# as trailing zeros survives. A bare `.Value = "2.60"` assignment gets
# auto-converted by Excel into the number 2.6, so we briefly force the
# cell to Text format for the assignment, then clear the format again so
# the cell is left with no explicit style (matching the original file).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.615.99"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.29%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.796.50"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.83%  "
# Row 4
$ws.Range("E4").Value = "  -0.06%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.89"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.39%  "
# Row 6
$ws.Range("E6").Value = "  +2.30%  "
# Row 7
$ws.Range("E7").Value = "  -0.09%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "33.03"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.23%  "
# Row 9
$ws.Range("E9").Value = "  +1.86%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0695"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.93%  "
# Row 11
$ws.Range("E11").Value = "  +0.43%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.055.53"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.78%  "
# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.804.26"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.15%  "
# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.09"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.79%  "
# Row 15
$ws.Range("E15").Value = "  +2.40%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.577.12"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.31%  "
# Row 17
$ws.Range("E17").Value = "  +3.07%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.89"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.46%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.72"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.65%  "
# Row 20
$ws.Range("E20").Value = "  +1.16%  "
# Row 21
$ws.Range("E21").Value = "  +2.91%  "
# Row 22
$ws.Range("E22").Value = "  -0.16%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.18"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.06%  "
# Row 24
$ws.Range("E24").Value = "  +1.69%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.63"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.79%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.32"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.93%  "
# Row 27
$ws.Range("E27").Value = "  +1.75%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.116"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.30%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.86%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.09"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +11.62%  "
# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.82"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.47%  "
# Row 32
$ws.Range("E32").Value = "  +1.21%  "
# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0525"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.10%  "
# Row 34
$ws.Range("E34").Value = "  +2.70%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.427.99"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.13%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.60"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +7.71%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.673"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.02%  "
# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.07"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.15%  "
# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0193"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.88%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "85.79"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +6.78%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.41"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.58%  "
# Row 42
$ws.Range("E42").Value = "  +1.20%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.75"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.02%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.68"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.81%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0529"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.99%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.12"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.00%  "
# Row 47
$ws.Range("E47").Value = "  +0.45%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.954.85"
$ws.Range("D48").ClearFormats()
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.13"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.32%  "
# Row 50
$ws.Range("E50").Value = "  -0.07%  "
# Row 51
$ws.Range("E51").Value = "  -5.08%  "
